# Fruta / hortaliza, semanal
# Insert 5 new weekly rows of "Macroferia Regional de Talca" - "Limón" price data
# just above row 1016, pushing the existing rows 1016-1028 down to 1021-1033.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows before row 1016 (shifts old rows 1016-1028 down to 1021-1033)
$insertRange = $ws.Range("A1016:T1020")
$insertRange.Insert()

# Fixed columns shared by every row in this block
$mercadoId = 5
$mercado = "Macroferia Regional de Talca"
$region = "Maule"
$periodicidad = 7
$rubro = "Fruta"
$grupoId = 100102
$grupo = "Cítricos"
$especieId = 100102003
$especie = "Limón"
$variedad = "Sin especificar"

function Set-LimonRow {
    param($row, $fecha, $calidad, $volumen, $precioMin, $precioProm, $precioMax, $unidad, $origen, $precioKilo, $kilos)

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = $periodicidad
    $ws.Cells.Item($row, 6).Value = $rubro
    $ws.Cells.Item($row, 7).Value = $grupoId
    $ws.Cells.Item($row, 8).Value = $grupo
    $ws.Cells.Item($row, 9).Value = $especieId
    $ws.Cells.Item($row, 10).Value = $especie
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precioMin
    $ws.Cells.Item($row, 15).Value = $precioProm
    $ws.Cells.Item($row, 16).Value = $precioMax
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $precioKilo
    $ws.Cells.Item($row, 20).Value = $kilos
}

Set-LimonRow 1016 44628 "1a plateado"  320 22000 22000 22000 "$/malla 14 kilos" "Provincia de Quillota"      1571 14
Set-LimonRow 1017 44628 "1a plateado"  360 22000 22000 22000 "$/malla 14 kilos" "Región de O'Higgins"        1571 14
Set-LimonRow 1018 44628 "2a amarillo"  190 19000 19000 19000 "$/malla 14 kilos" "Provincia de Quillota"      1357 14
Set-LimonRow 1019 44628 "2a amarillo"  200 19000 19000 19000 "$/malla 14 kilos" "Región de O'Higgins"        1357 14
Set-LimonRow 1020 44628 "3a amarillo"  100 12000 12000 12000 "$/malla 14 kilos" "Provincia de Quillota"       857 14
